$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 100
$ws.Range("I5").Value = 100
$ws.Range("K5").Value = 100
$ws.Range("M5").Value = 15

$ws.Range("H13").Value = 3002.5
$ws.Range("I13").Value = 3002.5
$ws.Range("K13").Value = 3002.5
$ws.Range("M13").Value = -2833.5

$ws.Range("H53").Value = 524.1
$ws.Range("I53").Value = 579
$ws.Range("K53").Value = 579
$ws.Range("M53").Value = 58

$ws.Range("H64").Value = 2583

$ws.Range("H67").Value = 2583

$ws.Range("H74").Value = 2249.75
$ws.Range("I74").Value = 2249.75
$ws.Range("K74").Value = 2249.75
$ws.Range("M74").Value = -1313.75

$ws.Range("H77").Value = 2249.75
$ws.Range("I77").Value = 2249.75
$ws.Range("K77").Value = 11248.75
$ws.Range("M77").Value = -6568.75

$ws.Range("H113").Value = 1999.5714
$ws.Range("I113").Value = 1999.6666
$ws.Range("K113").Value = 1999.6666
$ws.Range("M113").Value = 1254.3334

$ws.Range("H125").Value = 900
$ws.Range("I125").Value = 900
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 8100
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -5640
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3133.625
$ws.Range("I61").Value = 2678.1667
$ws.Range("K61").Value = 2678.1667
$ws.Range("M61").Value = -2466.1667

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()

$ws.Range("H110").Value = 867.5
$ws.Range("I110").Value = 867.5
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 867.5
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1177.5
$ws.Range("N110").ClearContents()

$ws.Range("H132").Value = 7013.5
$ws.Range("I132").Value = 5221.6
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 15664.8
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -13134.8
$ws.Range("N132").Value = -35060

$ws.Range("H136").Value = 3133.625
$ws.Range("I136").Value = 2678.1667
$ws.Range("K136").Value = 8034.500100000001
$ws.Range("M136").Value = -5484.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 20976.4
$ws.Range("I54").Value = 20976.4
$ws.Range("K54").Value = 20976.4
$ws.Range("M54").Value = -20492.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 12.1
$ws.Range("I7").Value = 2.3333333
$ws.Range("K7").Value = 2.3333333
$ws.Range("M7").Value = 110.6666667

$ws.Range("H31").Value = 1470
$ws.Range("I31").Value = 1470
$ws.Range("K31").Value = 1470
$ws.Range("M31").Value = -1175

$ws.Range("H34").Value = 1470
$ws.Range("I34").Value = 1470
$ws.Range("K34").Value = 1470
$ws.Range("M34").Value = -1268

$ws.Range("H43").Value = 5899
$ws.Range("J43").Value = 5899
$ws.Range("L43").Value = 5899
$ws.Range("N43").Value = -6267

$ws.Range("H58").Value = 4003
$ws.Range("I58").Value = 2575.7144
$ws.Range("J58").Value = 7333.3335
$ws.Range("K58").Value = 2575.7144
$ws.Range("L58").Value = 7333.3335
$ws.Range("M58").Value = -2372.7144
$ws.Range("N58").Value = -7739.3335

$ws.Range("H101").Value = 5899
$ws.Range("J101").Value = 5899
$ws.Range("L101").Value = 5899
$ws.Range("N101").Value = -12389

$ws.Range("H134").Value = 5370.6665
$ws.Range("I134").Value = 112
$ws.Range("K134").Value = 336
$ws.Range("M134").Value = 2199

$ws.Range("H136").Value = 4003
$ws.Range("I136").Value = 2575.7144
$ws.Range("J136").Value = 7333.3335
$ws.Range("K136").Value = 7727.1432
$ws.Range("L136").Value = 22000.0005
$ws.Range("M136").Value = -5177.1432
$ws.Range("N136").Value = -27100.0005

$ws.Range("H141").Value = 84999
$ws.Range("J141").Value = 84999
$ws.Range("L141").Value = 84999
$ws.Range("N141").Value = -95359

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 113.9
$ws.Range("I12").Value = 204.33333
$ws.Range("K12").Value = 612.99999
$ws.Range("M12").Value = -439.99999

$ws.Range("H38").Value = 1636.6
$ws.Range("I38").Value = 2647.6667
$ws.Range("J38").Value = 120
$ws.Range("K38").Value = 7943.000100000001
$ws.Range("L38").Value = 360
$ws.Range("M38").Value = -7596.000100000001
$ws.Range("N38").Value = -1054

$ws.Range("H68").Value = 801
$ws.Range("I68").Value = 801
$ws.Range("K68").Value = 2403
$ws.Range("M68").Value = -1592

$ws.Range("H71").Value = 801
$ws.Range("I71").Value = 801
$ws.Range("K71").Value = 7209
$ws.Range("M71").Value = -3153

$ws.Range("H131").Value = 4569.7144
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 4569.7144
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 13709.1432
$ws.Range("N131").Value = -23789.1432
$ws.Range("M131").ClearContents()

$ws.Range("H132").Value = 2400
$ws.Range("I132").Value = 2400
$ws.Range("K132").Value = 21600
$ws.Range("M132").Value = -19070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11748.25
$ws.Range("I40").Value = 8997.666999999999
$ws.Range("K40").Value = 8997.666999999999
$ws.Range("M40").Value = -8861.666999999999

$ws.Range("H55").Value = 877.875
$ws.Range("I55").Value = 590
$ws.Range("J55").Value = 1050.6
$ws.Range("K55").Value = 590
$ws.Range("L55").Value = 1050.6
$ws.Range("M55").Value = -417
$ws.Range("N55").Value = -1396.6

$ws.Range("H61").Value = 4366.25
$ws.Range("I61").Value = 3736.2
$ws.Range("J61").Value = 5416.3335
$ws.Range("K61").Value = 3736.2
$ws.Range("L61").Value = 5416.3335
$ws.Range("M61").Value = -3534.2
$ws.Range("N61").Value = -5820.3335

$ws.Range("H113").Value = 4366.25
$ws.Range("I113").Value = 3736.2
$ws.Range("J113").Value = 5416.3335
$ws.Range("K113").Value = 3736.2
$ws.Range("L113").Value = 5416.3335
$ws.Range("M113").Value = -1566.2
$ws.Range("N113").Value = -9756.333500000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 17666.334
$ws.Range("J101").Value = 17666.334
$ws.Range("L101").Value = 17666.334
$ws.Range("N101").Value = -24156.334

$ws.Range("H132").Value = 1937.5
$ws.Range("I132").Value = 1937.5
$ws.Range("K132").Value = 5812.5
$ws.Range("M132").Value = -3282.5
